# Update the "想去人数" (F column) counts for both the "展览" and "全部类型"
# worksheets. Both sheets contain the same rows of exhibition events, and
# the F-column values for several rows were incremented.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    3  = 271
    6  = 277
    7  = 6664
    11 = 80
    13 = 11
    16 = 217
    17 = 559
    18 = 61
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
